$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B50: convert from inline string "5" to numeric 5
$ws.Range("B50").Value = 5

# Add new row 51
$ws.Range("A51").Value = "Sunsi Wu"
$ws.Range("B51").Value = "'4"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "would be"
$ws.Range("D51").Value = "SUG"
$ws.Range("E51").Value = "WRI"
$ws.Range("F51").Value = "2a7301cf-d5b3-4d65-86b5-7931ca3b6163"
$ws.Range("G51").Value = "r1q7n9gAb_annotated.xlsx"
$ws.Range("H51").Value = "It would be beneficial for the clarity if authors define what they mean by convergence (normalised weight vector, angle, whichever path seems most natural) as early in the paper as possible."
